# Updates the "cryptos" price/volume table to the latest scrape.
# Note: several Price-column values look numeric (e.g. "0.0852", "1.00"),
# so they're written with a leading "'" (Excel's text-prefix) to keep the
# literal digits/trailing zeros as text instead of Excel coercing them to
# a float and dropping formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.287.15"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "3.585.98"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'606.09"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'147.87"
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("D7").Value = "3.583.73"
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.488"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "4.195.39"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").Value = "'29.51"
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("D16").Value = "3.581.60"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "'0.117"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "66.377.90"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "'11.03"
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("D20").Value = "'6.32"
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("D21").Value = "'14.84"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").Value = "'422.95"
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("D23").Value = "'0.609"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").Value = "'78.44"
$ws.Range("E26").Value = "  +3.18%  "
$ws.Range("D27").Value = "'8.22"
$ws.Range("E27").Value = "  +5.11%  "
$ws.Range("E28").Value = "  +3.11%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "3.582.74"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("E32").Value = "  +3.89%  "
$ws.Range("D33").Value = "'25.01"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'7.74"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").Value = "'5.57"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").Value = "'174.92"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("D40").Value = "'0.0852"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'0.880"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").Value = "'45.83"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("E44").Value = "  -3.82%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  +5.13%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'23.61"
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'7.13"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "'24.17"
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("E50").Value = "  -5.32%  "
$ws.Range("D51").Value = "'0.954"
$ws.Range("E51").Value = "  +2.70%  "
